$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '86.073.43'
$ws.Range("E2").Value = '  +6.66%  '
$ws.Range("D3").Value = '3.332.04'
$ws.Range("E3").Value = '  +3.20%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'219.12"
$ws.Range("E5").Value = '  +3.12%  '
$ws.Range("D6").Value = "'636.81"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = "'0.325"
$ws.Range("E7").Value = '  +14.81%  '
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("D10").Value = '3.334.16'
$ws.Range("E10").Value = '  +3.30%  '
$ws.Range("D11").Value = "'0.598"
$ws.Range("E11").Value = '  -3.13%  '
$ws.Range("E12").Value = '  +4.03%  '
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").Value = '3.928.85'
$ws.Range("E14").Value = '  +2.83%  '
$ws.Range("D15").Value = "'34.28"
$ws.Range("E15").Value = '  +4.24%  '
$ws.Range("D16").Value = "'5.41"
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '85.505.65'
$ws.Range("E17").Value = '  +6.05%  '
$ws.Range("D18").Value = '3.304.80'
$ws.Range("E18").Value = '  +2.43%  '
$ws.Range("D19").Value = "'14.70"
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("E20").Value = '  +5.61%  '
$ws.Range("D21").Value = "'441.36"
$ws.Range("D22").Value = "'9.19"
$ws.Range("E22").Value = '  -2.93%  '
$ws.Range("D23").Value = "'5.27"
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("D24").Value = "'7.38"
$ws.Range("E24").Value = '  +5.79%  '
$ws.Range("D25").Value = "'5.50"
$ws.Range("E25").Value = '  +13.85%  '
$ws.Range("D26").Value = "'12.28"
$ws.Range("E26").Value = '  +10.43%  '
$ws.Range("D27").Value = '3.482.46'
$ws.Range("E27").Value = '  +2.61%  '
$ws.Range("D28").Value = "'78.37"
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("E29").Value = '  +2.82%  '
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").Value = "'619.67"
$ws.Range("E31").Value = '  +9.34%  '
$ws.Range("D32").Value = "'0.168"
$ws.Range("E32").Value = '  +34.77%  '
$ws.Range("D33").Value = "'9.28"
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("E37").Value = '  -3.06%  '
$ws.Range("D38").Value = "'23.34"
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").Value = "'6.49"
$ws.Range("E39").Value = '  +10.39%  '
$ws.Range("D40").Value = "'0.420"
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").Value = "'21.27"
$ws.Range("E42").Value = '  +4.44%  '
$ws.Range("D43").Value = "'3.09"
$ws.Range("E43").Value = '  +11.32%  '
$ws.Range("D44").Value = "'2.05"
$ws.Range("E44").Value = '  +10.10%  '
$ws.Range("D45").Value = "'159.11"
$ws.Range("E45").Value = '  -4.11%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = "'190.16"
$ws.Range("E47").Value = '  -1.73%  '
$ws.Range("E48").Value = '  +1.21%  '
$ws.Range("D49").Value = "'45.24"
$ws.Range("E49").Value = '  +3.54%  '
$ws.Range("D50").Value = "'0.794"
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("D51").Value = "'26.64"
$ws.Range("E51").Value = '  +2.57%  '

Write-Host "Updated cryptos list"